# Add "Description" column (Hot Tomato sample menu item descriptions)
# to the Sidedish, Beverages and Dessert sheets.

$wb = $excel.ActiveWorkbook

# ---- Sidedish (sheet2): add column C "Description" ----
$ws = $wb.Worksheets.Item("Sidedish")
$ws.Range("C1").Value = "Description"
$ws.Range("C2").Value = "Freshly baked bread slices"
$ws.Range("C3").Value = "Traditional tasty fries"
$ws.Range("C4").Value = "Tasty mashed potato mixed with pepper gravy"
$ws.Range("C5").Value = "Nicely fried squid served with special sauce"
$ws.Range("C6").Value = "Chicken wings cooked with our special spice"
$ws.Columns.Item(3).ColumnWidth = 41
$ws.Range("C9").Select() | Out-Null

# ---- Beverages (sheet3): add column D "Description" ----
$ws = $wb.Worksheets.Item("Beverages")
$ws.Range("D1").Value = "Description"
$ws.Range("D2").Value = "Traditional Latte"
$ws.Range("D3").Value = "Specially brewed cappuccino"
$ws.Range("D4").Value = "Must-try signature espresso coffee"
$ws.Range("D5").Value = "Classic Australian coffee"
$ws.Range("D6").Value = "Lower suger, more healther"
$ws.Range("D7").Value = "Canned soda"
$ws.Range("D8").Value = "Canned soda"
$ws.Range("D9").Value = "Freshly juiced with great taste"
$ws.Columns.Item(4).ColumnWidth = 28.67
$ws.Range("D9").Select() | Out-Null

# ---- Dessert (sheet4): add column C "Description" ----
$ws = $wb.Worksheets.Item("Dessert")
$ws.Range("C1").Value = "Description"
$ws.Range("C2").Value = "Your perfect coffee mate"
$ws.Range("C3").Value = "A fruity iced blend of tropical mango and passion fruit juice with a hibiscus infusion."
$ws.Range("C4").Value = "Classical and unforgettable"
$ws.Range("C5").Value = "Scooping happiness to life"
$ws.Columns.Item(3).ColumnWidth = 31.83
$ws.Range("C25").Select() | Out-Null
